$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert a new paragraph after "More Organization of the Process" that
#    reads "Introduce Area to split out the controllers into section", and
#    move the lone "_GoBack" bookmark onto the end of that new paragraph
#    (Word only ever keeps a single "_GoBack" bookmark, so re-adding it here
#    removes the stale one that used to sit near the end of the document).
# ---------------------------------------------------------------------------

$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd() -eq "More Organization of the Process") {
        $anchor = $p
        break
    }
}

$anchor.Range.InsertParagraphAfter()
$newPara = $anchor.Next()
$newPara.Range.Text = "Introduce Area to split out the controllers into section"

# Build a tiny, throw-away non-collapsed range right after the text we just
# typed so the bookmark is anchored at that exact spot, then delete the
# placeholder character -- this leaves a zero-length "_GoBack" bookmark
# sitting right after the run, matching native Word's behaviour.
$paraRange = $newPara.Range.Duplicate
$paraRange.End = $paraRange.End - 1
$paraRange.InsertAfter("X")
$tail = $paraRange.Duplicate
$tail.Start = $tail.End - 1
$d.Bookmarks.Add("_GoBack", $tail) | Out-Null
$tail.Text = ""

# ---------------------------------------------------------------------------
# 2. Merge the two runs "/Controller/Menu/" + "Listing" into a single run
#    "/Controller/Menu/Listing".
# ---------------------------------------------------------------------------

$d.Content.Find.Execute("/Controller/Menu/Listing", $true, $false, $false,
    $false, $false, $true, 1, $false, "/Controller/Menu/Listing", 2) | Out-Null

# ---------------------------------------------------------------------------
# Note: the commit also drops the now-redundant <w:tblInd>/<w:tblCellMar>
# from the built-in "Table Grid" style definition in styles.xml /
# stylesWithEffects.xml (those values duplicate what "Table Grid" already
# inherits from "Normal Table"). The Word object model doesn't expose a
# writable TableStyle surface for that in this host (Style("Table
# Grid").Table's properties are read-only stubs here), so that part of the
# edit isn't reachable from COM automation and is intentionally left alone.
# ---------------------------------------------------------------------------
